$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows right after the header (before the old row 2),
# shifting the existing data rows (old 2-21) down to 7-26.
$ws.Rows.Item(2).Resize(5).Insert()

# The insert copies formatting from the row above; clear it so the new
# rows have no cell style, matching the rest of the data rows.
$ws.Range("A2:C6").ClearFormats()

# New data for the 5 newly inserted rows (rows 2-6)
$topData = @(
    @(2.571562051773072, -5.523353099822998, -4.892651081085205),
    @(2.563363254070282, -5.564052700996399, -4.925167679786682),
    @(2.571200489997864, -5.45090651512146, -4.94497549533844),
    @(2.582025349140167, -5.429405391216278, -4.891633093357086),
    @(2.521161556243896, -5.436496257781982, -4.74793529510498)
)

$r = 2
foreach ($row in $topData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# New data to append at the bottom (rows 27-31)
$bottomData = @(
    @(-2.611050009727472, -2.573673054575919, -8.294337868690496),
    @(0.61596310138703, -2.870795279741297, -8.901223957538587),
    @(2.079445004463198, -5.368536770343783, -3.524431616067886),
    @(0.3216586112976074, -3.676267147064209, -3.865855693817138),
    @(1.706132471561434, -4.47040206193924, -5.197765350341799)
)

$r = 27
foreach ($row in $bottomData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
